# Apply updated crypto price/volume data as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.412.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "'2.244.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'306.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'93.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.80%  "
$ws.Range("D7").Value = "'0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'34.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'2.395.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.837"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'13.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'44.078.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "'0.0₃0962"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'12.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'65.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "'3.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.91%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'236.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'38.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.28%  "
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "'153.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").Value = "'0.0798"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.97%  "
$ws.Range("D35").Value = "'0.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'0.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "'14.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.18%  "
$ws.Range("D40").Value = "'3.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "'1.736.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'80.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.18%  "
$ws.Range("D46").Value = "'99.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'4.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("D48").Value = "'56.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'69.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "
